# Update DM integration test fixture
#
# - Make header row (row 1) bold on every data sheet.
# - Resize columns (bestFit-style widths recomputed by Excel once headers turn bold).
# - Replace the generated UUIDs in column A (rows 2..N) on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CodeSchemes
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CodeSchemes")

$ws.Range("A1:N1").Font.Bold = $true

$ws.Columns.Item(1).ColumnWidth  = 32.285714285714285
$ws.Columns.Item(2).ColumnWidth  = 17.428571428571427
$ws.Columns.Item(3).ColumnWidth  = 25.714285714285715
$ws.Columns.Item(4).ColumnWidth  = 22.428571428571427
$ws.Columns.Item(5).ColumnWidth  = 14.142857142857142
$ws.Columns.Item(6).ColumnWidth  = 19.142857142857142
$ws.Columns.Item(7).ColumnWidth  = 19.142857142857142
$ws.Columns.Item(8).ColumnWidth  = 19.142857142857142
$ws.Columns.Item(9).ColumnWidth  = 20.714285714285715
$ws.Columns.Item(10).ColumnWidth = 24.0
$ws.Columns.Item(11).ColumnWidth = 19.142857142857142
$ws.Columns.Item(12).ColumnWidth = 15.714285714285714
$ws.Columns.Item(13).ColumnWidth = 20.714285714285715
$ws.Columns.Item(14).ColumnWidth = 27.285714285714285

$ws.Range("A2").Value = "45754185-b964-4112-8fc4-f6d69a1fe881"

# ---------------------------------------------------------------------------
# Codes
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Codes")

$ws.Range("A1:J1").Font.Bold = $true

$ws.Columns.Item(1).ColumnWidth  = 33.42857142857143
$ws.Columns.Item(2).ColumnWidth  = 17.428571428571427
$ws.Columns.Item(3).ColumnWidth  = 15.714285714285714
$ws.Columns.Item(4).ColumnWidth  = 14.142857142857142
$ws.Columns.Item(5).ColumnWidth  = 16.857142857142858
$ws.Columns.Item(6).ColumnWidth  = 19.142857142857142
$ws.Columns.Item(7).ColumnWidth  = 20.714285714285715
$ws.Columns.Item(8).ColumnWidth  = 24.0
$ws.Columns.Item(9).ColumnWidth  = 19.142857142857142
$ws.Columns.Item(10).ColumnWidth = 15.714285714285714

$ws.Range("A2").Value = "9ccf955b-cc12-419c-b039-0ef4e279fef3"
$ws.Range("A3").Value = "08172017-95dc-4e6c-b07e-fae9c695c8a8"
$ws.Range("A4").Value = "8627d6d9-4270-48ac-ab09-6796331e763c"
$ws.Range("A5").Value = "e9f3e1a7-38eb-42c8-bea8-dc7b359730a7"

# ---------------------------------------------------------------------------
# Extensions
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Extensions")

$ws.Range("A1:I1").Font.Bold = $true

$ws.Columns.Item(1).ColumnWidth = 31.142857142857142
$ws.Columns.Item(2).ColumnWidth = 17.428571428571427
$ws.Columns.Item(3).ColumnWidth = 14.142857142857142
$ws.Columns.Item(4).ColumnWidth = 24.0
$ws.Columns.Item(5).ColumnWidth = 15.714285714285714
$ws.Columns.Item(6).ColumnWidth = 19.142857142857142
$ws.Columns.Item(7).ColumnWidth = 19.142857142857142
$ws.Columns.Item(8).ColumnWidth = 15.714285714285714
$ws.Columns.Item(9).ColumnWidth = 24.0

$ws.Range("A2").Value = "4bd7767a-8afd-4317-b401-dc4dba7c4eb8"

# ---------------------------------------------------------------------------
# Members_dpmDimension
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Members_dpmDimension")

$ws.Range("A1:C1").Font.Bold = $true

$ws.Columns.Item(1).ColumnWidth = 33.42857142857143
$ws.Columns.Item(2).ColumnWidth = 10.857142857142858
$ws.Columns.Item(3).ColumnWidth = 30.571428571428573

$ws.Range("A2").Value = "c2d9797f-ae43-4531-ac94-321a70a2739c"
$ws.Range("A3").Value = "1ca3c132-e5f3-428b-85db-2379580931b2"
$ws.Range("A4").Value = "03efa971-0985-465b-ae41-1322c7bf6e87"
$ws.Range("A5").Value = "0571eb0a-fb93-4d12-a97a-1e3fabac269a"
